# Rename header cells on existing sheets
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Range("B1").Value = "Weekly_PO_Qty"
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" worksheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "PO Forecast"

# Copy header formatting (bold, centered, bordered) from sheet1's header row
$ws1.Range("A1:B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# Copy date-formatted cell style from sheet1's date column onto column A data rows
$ws1.Range("A2").Copy()
$ws3.Range("A2:A23").PasteSpecial(-4122)

# Populate forecast data
$ws3.Range("A2").Value = 45424.99999999999
$ws3.Range("B2").Value = 99
$ws3.Range("C2").Value = -127.3428001494996
$ws3.Range("D2").Value = 334.9201899103673
$ws3.Range("A3").Value = 45431.99999999999
$ws3.Range("B3").Value = 104
$ws3.Range("C3").Value = -123.434490089254
$ws3.Range("D3").Value = 363.3806657282719
$ws3.Range("A4").Value = 45445.99999999999
$ws3.Range("B4").Value = 116
$ws3.Range("C4").Value = -117.9715691109359
$ws3.Range("D4").Value = 348.4687526312036
$ws3.Range("A5").Value = 45452.99999999999
$ws3.Range("B5").Value = 122
$ws3.Range("C5").Value = -122.9410330464104
$ws3.Range("D5").Value = 372.4807312388508
$ws3.Range("A6").Value = 45459.99999999999
$ws3.Range("B6").Value = 127
$ws3.Range("C6").Value = -106.0614382893155
$ws3.Range("D6").Value = 365.3147158610375
$ws3.Range("A7").Value = 45466.99999999999
$ws3.Range("B7").Value = 133
$ws3.Range("C7").Value = -85.1111986912761
$ws3.Range("D7").Value = 387.8765943828211
$ws3.Range("A8").Value = 45529.99999999999
$ws3.Range("B8").Value = 185
$ws3.Range("C8").Value = -55.47389235347871
$ws3.Range("D8").Value = 422.8973932271363
$ws3.Range("A9").Value = 45550.99999999999
$ws3.Range("B9").Value = 202
$ws3.Range("C9").Value = -22.19911263794116
$ws3.Range("D9").Value = 428.6410866486146
$ws3.Range("A10").Value = 45557.99999999999
$ws3.Range("B10").Value = 207
$ws3.Range("C10").Value = -24.65590321279046
$ws3.Range("D10").Value = 443.8989127549382
$ws3.Range("A11").Value = 45564.99999999999
$ws3.Range("B11").Value = 213
$ws3.Range("C11").Value = -10.43303189965272
$ws3.Range("D11").Value = 450.813325808946
$ws3.Range("A12").Value = 45571.99999999999
$ws3.Range("B12").Value = 219
$ws3.Range("C12").Value = -14.68418044822339
$ws3.Range("D12").Value = 469.3489422657867
$ws3.Range("A13").Value = 45578.99999999999
$ws3.Range("B13").Value = 225
$ws3.Range("C13").Value = -7.570189217888969
$ws3.Range("D13").Value = 458.9419160700058
$ws3.Range("A14").Value = 45585.99999999999
$ws3.Range("B14").Value = 230
$ws3.Range("C14").Value = -5.320195432278322
$ws3.Range("D14").Value = 464.9863102379774
$ws3.Range("A15").Value = 45592.99999999999
$ws3.Range("B15").Value = 236
$ws3.Range("C15").Value = -10.43579987481106
$ws3.Range("D15").Value = 465.2050951774054
$ws3.Range("A16").Value = 45599.99999999999
$ws3.Range("B16").Value = 242
$ws3.Range("C16").Value = -2.521435943177874
$ws3.Range("D16").Value = 476.4655427390982
$ws3.Range("A17").Value = 45606.99999999999
$ws3.Range("B17").Value = 247
$ws3.Range("C17").Value = 16.37155562017537
$ws3.Range("D17").Value = 483.3903579090415
$ws3.Range("A18").Value = 45613.99999999999
$ws3.Range("B18").Value = 253
$ws3.Range("C18").Value = 23.25898669137521
$ws3.Range("D18").Value = 480.3972292059281
$ws3.Range("A19").Value = 45620.99999999999
$ws3.Range("B19").Value = 259
$ws3.Range("C19").Value = 19.55134733413933
$ws3.Range("D19").Value = 490.1622139695514
$ws3.Range("A20").Value = 45627.99999999999
$ws3.Range("B20").Value = 265
$ws3.Range("C20").Value = 34.60993271631352
$ws3.Range("D20").Value = 512.1994345425774
$ws3.Range("A21").Value = 45634.99999999999
$ws3.Range("B21").Value = 270
$ws3.Range("C21").Value = 15.62722474092234
$ws3.Range("D21").Value = 496.5629107741404
$ws3.Range("A22").Value = 45641.99999999999
$ws3.Range("B22").Value = 276
$ws3.Range("C22").Value = 52.93021460266859
$ws3.Range("D22").Value = 509.4950351772831
$ws3.Range("A23").Value = 45648.99999999999
$ws3.Range("B23").Value = 282
$ws3.Range("C23").Value = 48.56006347362266
$ws3.Range("D23").Value = 526.4085266717534

$excel.CutCopyMode = $false
